$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 2 (pushes existing rows 2-10 down to 5-13)
$ws.Rows.Item(2).Resize(3).Insert()

# Insert a column before E (pushes old E/F/G -> F/G/H)
$ws.Columns.Item(5).Insert()
# Insert a column before the new G (old F/attributes, now at G) -> pushes it to H, old G(degree) -> I
$ws.Columns.Item(7).Insert()

# New "Skill class" style list of strings in column H, rows 2-10
$ws.Range("H2").Value = "scale"
$ws.Range("H3").Value = "shield"
$ws.Range("H4").Value = "uhp"
$ws.Range("H5").Value = "mana"
$ws.Range("H6").Value = "speed"
$ws.Range("H7").Value = "armor"
$ws.Range("H8").Value = "damage"
$ws.Range("H9").Value = "freq"
$ws.Range("H10").Value = "range"

# New row 4 ("damage" attribute/value pair), matches style of rows 5-10 (C/D pairs)
$ws.Range("C4").Value = "damage"
$ws.Range("D4").Value = -1

# New header cells (row 11) for the inserted columns E and G
$ws.Range("E11").Value = "tar_type"
$ws.Range("G11").Value = "cold_t"

# New data cells for row 12 (Marine / Shield Defense)
$ws.Range("E12").Value = "n"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0

# New data cells for row 13 (Rockhead / Shocking)
$ws.Range("E13").Value = "n"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0

# Match the new active-cell selection recorded in the file
$ws.Range("H13").Select()
